# Add season-record columns (Wins / Losses / Ties) to the sheet, one
# column each, with the team's full-season record repeated on every
# player row - this is the "season record" data the commit message
# says was missing from the original scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ------------------------------------------------
# Give the three new headers the same look as the rest of the header
# row (bold font, thin border, centered/top aligned) by copying the
# formatting from the last existing header cell (AC1).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows (2-48) --------------------------------------------------
# Every player belongs to the same team/season, so the record is
# identical (90 wins, 72 losses, 0 ties) on each row.
for ($row = 2; $row -le 48; $row++) {
    $ws.Cells.Item($row, 30).Value = 90
    $ws.Cells.Item($row, 31).Value = 72
    $ws.Cells.Item($row, 32).Value = 0
}
